$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the token name/type/P-T rows for each token into a single
# Python-tuple-style string in column A, collapsing the previous
# multi-row-per-token layout down to one row per token.
$ws.Range("A2").Value = "('Faerie Rogue', ['Token Creature — Faerie Rogue', 'Flying', '1/1'])"
$ws.Range("A3").Value = "('Giant Warrior', ['Token Creature — Giant Warrior', '5/5'])"
$ws.Range("A4").Value = "('Treefolk Shaman', ['Token Creature — Treefolk Shaman', '2/5'])"

# The old rows 5-11 (which held the now-merged detail lines) are no
# longer needed; delete them and shift the remaining rows up so the
# used range shrinks to A1:A4.
$ws.Range("A5:A11").EntireRow.Delete()
